# Adds a new "2021年" data row (row 5) to Sheet1, mirroring the structure
# of the existing yearly rows (2018年-2020年) in rows 2-4.
#
# NOTE on blank cells: a handful of columns (G, L, R, S, AD, AG, AS, AT,
# BC, CJ, CU, CW, CZ, DK) have no data point for 2021 and are left
# untouched (i.e. not written), matching how the existing sheet already
# represents "no value" cells for these indicator columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @("A", "S", "2021年"),
    @("B", "N", 52.9),
    @("C", "N", 103.8),
    @("D", "N", -6.2),
    @("E", "N", 93.09999999999999),
    @("F", "N", -2),
    @("G", "E", $null),
    @("H", "N", 179.7),
    @("I", "N", 149.6),
    @("J", "N", 149.5),
    @("K", "N", -19.5),
    @("L", "E", $null),
    @("M", "N", 28.7),
    @("N", "N", -70.5),
    @("O", "N", -4.5),
    @("P", "N", 63.7),
    @("Q", "N", -81),
    @("R", "E", $null),
    @("S", "E", $null),
    @("T", "N", 2.9),
    @("U", "N", 3),
    @("V", "N", 38.3),
    @("W", "N", -1.9),
    @("X", "N", 30.4),
    @("Y", "N", 8),
    @("Z", "N", 171.6),
    @("AA", "N", -34.2),
    @("AB", "N", 21.1),
    @("AC", "N", 22.2),
    @("AD", "E", $null),
    @("AE", "N", 10.3),
    @("AF", "N", -71.3),
    @("AG", "E", $null),
    @("AH", "N", -14.8),
    @("AI", "N", -24.2),
    @("AJ", "N", -34.1),
    @("AK", "N", 125.9),
    @("AL", "N", -16.9),
    @("AM", "N", 55),
    @("AN", "N", -18),
    @("AO", "N", -4.8),
    @("AP", "N", -46.3),
    @("AQ", "N", -36.8),
    @("AR", "N", 5.2),
    @("AS", "E", $null),
    @("AT", "E", $null),
    @("AU", "N", 1.9),
    @("AV", "N", 3.3),
    @("AW", "N", 157.8),
    @("AX", "N", 22.3),
    @("AY", "N", -9.199999999999999),
    @("AZ", "N", 6.6),
    @("BA", "N", -12.1),
    @("BB", "N", -3.4),
    @("BC", "E", $null),
    @("BD", "N", -87.8),
    @("BE", "N", 14.4),
    @("BF", "N", 203.1),
    @("BG", "N", 46.7),
    @("BH", "N", 9.4),
    @("BI", "N", -35.6),
    @("BJ", "N", 1),
    @("BK", "N", -32.5),
    @("BL", "N", -5.2),
    @("BM", "N", -5.8),
    @("BN", "N", -18),
    @("BO", "N", 144.8),
    @("BP", "N", -37.4),
    @("BQ", "N", -25.6),
    @("BR", "N", 151.2),
    @("BS", "N", -41.6),
    @("BT", "N", -9.4),
    @("BU", "N", -65.8),
    @("BV", "N", -20.2),
    @("BW", "N", -19.9),
    @("BX", "N", -4.4),
    @("BY", "N", -29.3),
    @("BZ", "N", 59.7),
    @("CA", "N", 53),
    @("CB", "N", -38.1),
    @("CC", "N", 73.7),
    @("CD", "N", -71.7),
    @("CE", "N", -0.4),
    @("CF", "N", 66.3),
    @("CG", "N", 65.7),
    @("CH", "N", 28),
    @("CI", "N", 10.3),
    @("CJ", "E", $null),
    @("CK", "N", 136.1),
    @("CL", "N", 34.3),
    @("CM", "N", 19.5),
    @("CN", "N", -15.6),
    @("CO", "N", 91.7),
    @("CP", "N", 53),
    @("CQ", "N", -89.59999999999999),
    @("CR", "N", -38.1),
    @("CS", "N", 133.3),
    @("CT", "N", 10.1),
    @("CU", "E", $null),
    @("CV", "N", -6),
    @("CW", "E", $null),
    @("CX", "N", -49.2),
    @("CY", "N", 27.6),
    @("CZ", "E", $null),
    @("DA", "N", 131.5),
    @("DB", "N", 1.2),
    @("DC", "N", -17.2),
    @("DD", "N", 24.3),
    @("DE", "N", -27.4),
    @("DF", "N", -30.9),
    @("DG", "N", 85.59999999999999),
    @("DH", "N", -34.5),
    @("DI", "N", 148.6),
    @("DJ", "N", -1.6),
    @("DK", "E", $null),
)

$targetRow = 5

foreach ($entry in $rowData) {
    $col = $entry[0]
    $kind = $entry[1]
    $value = $entry[2]
    $cellRef = "$col$targetRow"

    if ($kind -eq "E") {
        # No data point for this column in the new row - leave the cell
        # blank, same as the sheet's other "missing value" cells.
        continue
    } elseif ($kind -eq "S") {
        $ws.Range($cellRef).Value = $value
    } else {
        $ws.Range($cellRef).Value = [double]$value
    }
}

# Match the formatting (bold label, border, centered/top aligned) that the
# other year-label cells in column A (A2:A4) use, by copying A4's format
# onto the new A5 label cell.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
